$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# Create the new "RuleConfig_Update" sheet by copying the structure
# of "RuleConfig_Delete" (same headers / row2 layout), then place it
# right after "RuleConfig_Delete" (i.e. as the last tab).
# ------------------------------------------------------------------
$wsDelete = $wb.Worksheets.Item("RuleConfig_Delete")
$wsDelete.Copy([System.Type]::Missing, $wsDelete) | Out-Null

$wsUpdate = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsUpdate.Name = "RuleConfig_Update"

# New sheet's HTTP Type (C2) becomes PUT - set first so it lands at
# shared-string index 37 (the order in which brand-new strings are
# introduced controls their shared-string index, matching the diff).
$wsUpdate.Range("C2").Value = "PUT"
$wsUpdate.Range("C2").Select() | Out-Null

# ------------------------------------------------------------------
# RuleConfig_Delete : HTTP Type becomes DELETE
# ------------------------------------------------------------------
$wsDelete.Range("C2").Value = "DELETE"
$wsDelete.Range("C2").Select() | Out-Null

# ------------------------------------------------------------------
# RuleConfig_View : HTTP Type becomes GET
# ------------------------------------------------------------------
$wsView = $wb.Worksheets.Item("RuleConfig_View")
$wsView.Range("C2").Value = "GET"
$wsView.Range("C2").Select() | Out-Null

# Re-activate the new sheet so it ends up as the active tab, with C2
# selected (matches the saved selection state for RuleConfig_Update).
$wsUpdate.Activate() | Out-Null
$wsUpdate.Range("C2").Select() | Out-Null
